$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The activity stats for "Tushar Deshpande" rows 2 and 3 were swapped:
# row 2 now holds the 20-run innings (runs/balls/fours/sixes),
# row 3 now holds the 1-run innings.
# These columns are stored as text (number-stored-as-text), so force a
# text number format before writing the values back.
$ws.Range("C2:F3").NumberFormat = "@"

$ws.Range("C2").Value = "20"
$ws.Range("D2").Value = "9"
$ws.Range("E2").Value = "2"
$ws.Range("F2").Value = "1"

$ws.Range("C3").Value = "1"
$ws.Range("D3").Value = "3"
$ws.Range("E3").Value = "0"
$ws.Range("F3").Value = "0"
